# Update the "dSF" (column F) values for several rows as per the repull/push of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -4
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -3
$ws.Range("F12").Value = -3
$ws.Range("F13").Value = -12
$ws.Range("F16").Value = -1
